$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the date text in the "Ternate , <date>" cell
$ws.Range("E29").Value = "Ternate , 29 Agustus 2023"

# Fill previously-empty cells with explicit zero values
$ws.Range("E13").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("E16").Value = 0
$ws.Range("E17").Value = 0
